$d = $word.ActiveDocument

$pairs = @(
    @("152×5=760", "593×8=4744"),
    @("685×6=4110", "844×4=3376"),
    @("141×4=564", "257×8=2056"),
    @("120×6=720", "804×9=7236"),
    @("742×6=4452", "490×4=1960"),
    @("790×5=3950", "975×6=5850"),
    @("980×7=6860", "185×2=370"),
    @("153×5=765", "568×3=1704"),
    @("572×4=2288", "854×3=2562"),
    @("709×4=2836", "972×8=7776"),
    @("441×3=1323", "596×5=2980"),
    @("905×6=5430", "121×9=1089"),
    @("317×9=2853", "308×6=1848"),
    @("560×4=2240", "843×3=2529"),
    @("247×6=1482", "592×6=3552"),
    @("910×8=7280", "141×7=987"),
    @("481×6=2886", "928×7=6496"),
    @("779×4=3116", "707×9=6363"),
    @("968×8=7744", "225×9=2025"),
    @("706×4=2824", "935×4=3740"),
    @("808×5=4040", "943×3=2829"),
    @("175×8=1400", "506×5=2530"),
    @("256×6=1536", "263×2=526"),
    @("882×9=7938", "635×2=1270"),
    @("138×9=1242", "787×2=1574")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
